$wb = $excel.ActiveWorkbook

# --- Update Sample Name values on the "Samples" sheet (hyphens -> underscores) ---
$samples = $wb.Worksheets.Item("Samples")
$samples.Range("A2").Value = "His_neg_M3_T02_liv"
$samples.Range("A3").Value = "His_neg_M3_T03_kid"
$samples.Range("A4").Value = "His_neg_M3_T22_serum_30m"
$samples.Range("A5").Value = "His_neg_M3_T22_serum_120m"

# --- Widen / refresh the remembered selection on every sheet to A2:A5 ---
$animals = $wb.Worksheets.Item("Animals")
$animals.Activate()
$animals.Range("A2:A5").Select()

$treatments = $wb.Worksheets.Item("Treatments")
$treatments.Activate()
$treatments.Range("A2:A5").Select()

$tissues = $wb.Worksheets.Item("Tissues")
$tissues.Activate()
$tissues.Range("A2:A5").Select()

# --- Samples becomes the active / tab-selected sheet, selection A2:A5 ---
$samples.Activate()
$samples.Range("A2:A5").Select()
